$d = $word.ActiveDocument

# Locate the paragraph that reads "Quảng Trị, <ngay_thang>" immediately
# above the "Ý KIẾN LÃNH ĐẠO CỤC" block (the second such paragraph in the
# document; the first, near the letterhead, must stay untouched). Keep
# scanning to the end so `$target` ends up on the *last* match.
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Quảng Trị,*<ngay_thang>*") {
        $target = $p
    }
}

if ($target -ne $null) {
    # Delete the whole paragraph, including its paragraph mark, merging the
    # following paragraph up into its place.
    $d.Range($target.Range.Start, $target.Range.End).Delete()
}

# Move the "_GoBack" bookmark that currently sits just before "Thành phần
# kiểm tra gồm có:" so that it instead starts the paragraph that now
# follows directly after the deleted "Quảng Trị, <ngay_thang>" paragraph
# (the short, tab-only paragraph right before "Ý KIẾN LÃNH ĐẠO CỤC").
$dest = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Ý KIẾN L*ÃNH*ĐẠO*C*ỤC*") {
        $dest = $p
        break
    }
}

if ($dest -ne $null) {
    $destPara = $dest.Previous()
    $pt = $d.Range($destPara.Range.Start, $destPara.Range.Start)
    $d.Bookmarks.Add("_GoBack", $pt)
}
